$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the footprint and JLCPCB part number for the common-mode choke (L1)
# from the 1206 package variant to the 0805 package variant.
$ws.Range("C34").Value = "L_CommonMode_Wurth_WE-CNSW-0805"
$ws.Range("D34").Value = "C2649324"

# Move the active selection to D39 (matches the author's updated cursor
# position saved with the workbook).
$ws.Range("D39").Select()
